$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Analysis_Unit")
$ws2 = $wb.Worksheets.Item("r AnalysisUnit_Variable")

# --- Shift existing data rows 84..159 down by 3 rows (84..159 -> 87..162) ---
# Work bottom-up so we never overwrite a source row before it has been read.
# NOTE: the `.Value` getter in this COM shim does not resolve properly (it
# echoes a member signature instead of the cell's contents) - `.Value2` must
# be used for reads; `.Value` remains fine (and is used) for writes.
for ($r = 159; $r -ge 84; $r--) {
    $newR = $r + 3
    $bVal = $ws2.Cells.Item($r, 2).Value2
    $cVal = $ws2.Cells.Item($r, 3).Value2
    $fVal = $ws2.Cells.Item($r, 6).Value2

    $ws2.Cells.Item($newR, 2).Value = $bVal
    $ws2.Cells.Item($newR, 3).Value = $cVal
    $ws2.Cells.Item($newR, 6).Value = $fVal
}

# --- Write the 3 new variable rows into the freshly-vacated 84..86 slots ---
$newVars = "DM_Color_EDF", "DM_Color_CDS", "DM_Color_BOND"
for ($i = 0; $i -lt $newVars.Length; $i++) {
    $r = 84 + $i
    $name = $newVars[$i]
    $ws2.Cells.Item($r, 2).Value = $name
    $ws2.Cells.Item($r, 3).Value = $name
    $ws2.Cells.Item($r, 6).Value = $name
}

# --- Columns A ("Action") / E ("Main Analysis_Unit") are constant across the
#     whole 84..162 block (they also now cover rows 160..162, which are brand
#     new row positions pushed past the former end of the sheet at row 159). ---
for ($r = 84; $r -le 162; $r++) {
    $ws2.Cells.Item($r, 1).Value = "CREATE/MODIFY"
    $ws2.Cells.Item($r, 5).Value = "CUSTOMER"
}

# --- View/selection state: sheet2 becomes the active tab ---
$ws2.Activate()
$ws2.Range("F86").Select()

Write-Host "done"
